$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Output_flows")
$ws1.Range("C7").Value = 3.183813311787415 * [Math]::Pow(10, -24)
$ws1.Range("E7").Value = 1.984293364045812 * [Math]::Pow(10, -23)
$ws1.Range("F7").Value = 5.357592082923691 * [Math]::Pow(10, -24)
$ws1.Range("C12").Value = 1.341933390787962 * [Math]::Pow(10, -27)
$ws1.Range("E12").Value = 4.181761399881104 * [Math]::Pow(10, -27)
$ws1.Range("F12").Value = 1.129075577967898 * [Math]::Pow(10, -27)
$ws1.Range("C13").Value = 3.701176254331037 * [Math]::Pow(10, -28)
$ws1.Range("D13").Value = 2.742980459372967 * [Math]::Pow(10, -29)
$ws1.Range("E13").Value = 1.153368423557145 * [Math]::Pow(10, -27)
$ws1.Range("F13").Value = 3.114094743604292 * [Math]::Pow(10, -28)
$ws1.Range("C14").Value = 3.022253754609902 * [Math]::Pow(10, -28)
$ws1.Range("D14").Value = 8.959295556336483 * [Math]::Pow(10, -29)
$ws1.Range("E14").Value = 9.418011488820912 * [Math]::Pow(10, -28)
$ws1.Range("F14").Value = 2.542863101981646 * [Math]::Pow(10, -28)
$ws1.Range("C15").Value = 2.016736319143734 * [Math]::Pow(10, -27)
$ws1.Range("D15").Value = 5.081722938591634 * [Math]::Pow(10, -27)
$ws1.Range("E15").Value = 6.284596650644148 * [Math]::Pow(10, -27)
$ws1.Range("F15").Value = 1.69684109567392 * [Math]::Pow(10, -27)
$ws1.Range("C17").Value = 2.130540900317189 * [Math]::Pow(10, -19)
$ws1.Range("E17").Value = 2.655694763682091 * [Math]::Pow(10, -18)
$ws1.Range("F17").Value = 7.170375861941645 * [Math]::Pow(10, -19)
$ws1.Range("C18").Value = 4.867349315448403 * [Math]::Pow(10, -20)
$ws1.Range("D18").Value = 1.442915221480054 * [Math]::Pow(10, -22)
$ws1.Range("E18").Value = 6.067095021796351 * [Math]::Pow(10, -19)
$ws1.Range("F18").Value = 1.638115655885015 * [Math]::Pow(10, -19)
$ws1.Range("C19").Value = 3.500210998188149 * [Math]::Pow(10, -20)
$ws1.Range("D19").Value = 4.150469608151817 * [Math]::Pow(10, -22)
$ws1.Range("E19").Value = 4.362972810466508 * [Math]::Pow(10, -19)
$ws1.Range("F19").Value = 1.178002658825957 * [Math]::Pow(10, -19)
$ws1.Range("C20").Value = 8.906576523497707 * [Math]::Pow(10, -19)
$ws1.Range("D20").Value = 8.977029631164394 * [Math]::Pow(10, -20)
$ws1.Range("E20").Value = 1.110194534743045 * [Math]::Pow(10, -17)
$ws1.Range("F20").Value = 2.997525243806222 * [Math]::Pow(10, -18)

$ws2 = $wb.Worksheets.Item("Input_flows")
$ws2.Range("C7").Value = 2.838433903516922 * [Math]::Pow(10, -23)
$ws2.Range("C12").Value = 3.531510323106573 * [Math]::Pow(10, -27)
$ws2.Range("C13").Value = 8.011842624443985 * [Math]::Pow(10, -28)
$ws2.Range("C14").Value = 5.715612023862837 * [Math]::Pow(10, -28)
$ws2.Range("C15").Value = 1.507989700405344 * [Math]::Pow(10, -26)
$ws2.Range("C17").Value = 3.531572447118433 * [Math]::Pow(10, -18)
$ws2.Range("C18").Value = 8.011772697020322 * [Math]::Pow(10, -19)
$ws2.Range("C19").Value = 5.715606446096145 * [Math]::Pow(10, -19)
$ws2.Range("C20").Value = 1.507989853989808 * [Math]::Pow(10, -17)
